# Generate Report for handback
# Refresh the handoff/handback timestamps recorded on the per-locale report
# sheets ("zh-cn" and "de-de") for the 05bc5cc6... file, simulating a
# newly generated handback report with updated Correspond Handoff Datetime
# (column D) and Correspond Handback DateTime (column G) values.

$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D2").Value = "2016-01-14 05:04:30"   # Correspond Handoff Datetime
$wsZh.Range("G2").Value = "2016-01-14 05:05:56"   # Correspond Handback DateTime

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D2").Value = "2016-01-14 05:04:55"   # Correspond Handoff Datetime
$wsDe.Range("G2").Value = "2016-01-14 05:06:35"   # Correspond Handback DateTime
